$d = $word.ActiveDocument
$vt = [char]11

# Step 1: collapse the whole document to a single empty paragraph
$end = $d.Content.End
$d.Range(0, $end).Text = ""
while ($d.Paragraphs.Count -gt 1) {
    $p1end = $d.Paragraphs(1).Range.End
    $markR = $d.Range($p1end - 1, $p1end)
    $markR.Delete()
}

# Step 2: rebuild the target paragraphs token by token so whitespace runs get
#         proper xml:space handling from the COM host
function Append-Text([string]$text) {
    $e = $d.Content.End
    $r = $d.Range($e - 1, $e - 1)
    $r.InsertAfter($text)
}
function Append-Break() {
    $e = $d.Content.End
    $r = $d.Range($e - 1, $e - 1)
    $r.InsertAfter($vt)
}
function Append-Para() {
    $e = $d.Content.End
    $r = $d.Range($e - 1, $e - 1)
    $r.InsertParagraphAfter()
}

Append-Text 'Password Test Case SDET - Automated Test Code (Generated by Haiku)'

Append-Para
Append-Break
Append-Text 'Python Test Code (auto-generated from Test Lead requirements):'

Append-Para
Append-Text 'Here is the Python test function code:'
Append-Break
Append-Break
Append-Text '```python'
Append-Break
Append-Text 'import re'
Append-Break
Append-Break
Append-Text 'def test_password(password):'
Append-Break
Append-Text '  length_valid = len(password) >= 8'
Append-Break
Append-Text '  has_digit = re.search(r''\d'', password) is not None'
Append-Break
Append-Text '  has_special = re.search(r''[!@#$%^&*]'', password) is not None'
Append-Break
Append-Break
Append-Text '  if not length_valid:'
Append-Break
Append-Text '    print("Failed - Must be at least 8 characters")'
Append-Break
Append-Text '  elif not has_digit: '
Append-Break
Append-Text '    print("Failed - Must include at least one number")'
Append-Break
Append-Text '  elif not has_special:'
Append-Break
Append-Text '    print("Failed - Must include at least one special character")'
Append-Break
Append-Text '  else:'
Append-Break
Append-Text '    print("Passed")'
Append-Break
Append-Break
Append-Text 'test_password("abc1$")'
Append-Break
Append-Text 'test_password("abcd@xyz")  '
Append-Break
Append-Text 'test_password("abcd1234")'
Append-Break
Append-Text 'test_password("abcd@xyz")'
Append-Break
Append-Text 'test_password("abc1@def")'
Append-Break
Append-Text 'test_password("MyPass123!")'
Append-Break
Append-Text 'test_password("1234@5678")'
Append-Break
Append-Text 'test_password("abcdefgh")'
Append-Break
Append-Text 'test_password("abcd1234")'
Append-Break
Append-Text 'test_password("Ab1$xyz9")'
Append-Break
Append-Text '```'
Append-Break
Append-Break
Append-Text 'This implements test cases to validate the password requirements without any explanatory comments, as requested. It uses regex to check for digits and special characters.'

